# Auto-generated edit script: updates '想去人数' (want-to-go count, column F)
# values across the '展览', '演出' and '全部类型' worksheets, matching the
# upstream data refresh recorded in the commit diff.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 254
$ws.Range("F3").Value = 244
$ws.Range("F4").Value = 263
$ws.Range("F5").Value = 2867
$ws.Range("F6").Value = 64
$ws.Range("F8").Value = 2221
$ws.Range("F9").Value = 361
$ws.Range("F11").Value = 432
$ws.Range("F13").Value = 2548
$ws.Range("F15").Value = 1344
$ws.Range("F16").Value = 4694
$ws.Range("F18").Value = 5102
$ws.Range("F19").Value = 1665
$ws.Range("F20").Value = 2866
$ws.Range("F21").Value = 3266
$ws.Range("F22").Value = 166
$ws.Range("F23").Value = 1554
$ws.Range("F24").Value = 257
$ws.Range("F25").Value = 836
$ws.Range("F26").Value = 106
$ws.Range("F27").Value = 290
$ws.Range("F28").Value = 995
$ws.Range("F29").Value = 1857
$ws.Range("F30").Value = 117
$ws.Range("F31").Value = 280
$ws.Range("F32").Value = 704
$ws.Range("F33").Value = 156
$ws.Range("F34").Value = 336
$ws.Range("F35").Value = 414

# --- 演出 (sheet 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 98
$ws.Range("F8").Value = 90
$ws.Range("F10").Value = 21
$ws.Range("F11").Value = 194

# --- 全部类型 (sheet 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 98
$ws.Range("F7").Value = 254
$ws.Range("F8").Value = 244
$ws.Range("F10").Value = 263
$ws.Range("F11").Value = 2867
$ws.Range("F12").Value = 64
$ws.Range("F13").Value = 2221
$ws.Range("F14").Value = 361
$ws.Range("F15").Value = 90
$ws.Range("F18").Value = 432
$ws.Range("F20").Value = 21
$ws.Range("F21").Value = 2548
$ws.Range("F22").Value = 1344
$ws.Range("F23").Value = 194
$ws.Range("F26").Value = 4694
$ws.Range("F28").Value = 5102
$ws.Range("F29").Value = 1665
$ws.Range("F30").Value = 2866
$ws.Range("F31").Value = 3266
$ws.Range("F32").Value = 166
$ws.Range("F35").Value = 1554
$ws.Range("F37").Value = 257
$ws.Range("F38").Value = 836
$ws.Range("F39").Value = 106
$ws.Range("F40").Value = 290
$ws.Range("F41").Value = 995
$ws.Range("F43").Value = 1857
$ws.Range("F44").Value = 117
$ws.Range("F45").Value = 280
$ws.Range("F46").Value = 704
$ws.Range("F47").Value = 156
$ws.Range("F48").Value = 336
$ws.Range("F49").Value = 414

